# 1.1.1.xlsx — add the "2020" column (Q) to the single data sheet,
# mirroring the formatting of the existing "2019" column (P), and move
# the live selection to N19 (as captured in the saved sheetView).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New 2020 values, keyed by row -> value (column P holds 2019 today).
$newValues = @{
    4  = 2020
    5  = 0.02
    6  = 0
    7  = 0
    8  = 0
    9  = 0.54
    10 = 0
    11 = 0
    12 = 0
    13 = 0
    14 = 0
}

foreach ($row in 4..14) {
    $source = $ws.Cells.Item($row, 16)   # column P (16th column) = 2019
    $target = $ws.Cells.Item($row, 17)   # column Q (17th column) = 2020

    # Copy the 2019 cell's formatting (font, borders, alignment, number
    # format) onto the new 2020 cell, then overwrite with the real value.
    $source.Copy($target)
    $target.Value = $newValues[$row]
}

# Move / record the active selection, as shown in the saved sheetView.
$ws.Range("N19").Select()
